$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header/id values updated (meanEMG / legmaxROM column relabel)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) - B2 & D2 recomputed, C2 & E2 removed entirely
$ws.Range("B2").Value = 21.163384158590951
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 28.997406510130791
$ws.Range("E2").ClearContents()

# Row 3 (STR) - B3 recomputed, C3 removed, D3 newly added, E3 recomputed
$ws.Range("B3").Value = 19.524893642815496
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 21.237197820286276
$ws.Range("E3").Value = -10.160790160633542

# Selection now covers the updated data block only
$ws.Range("B1:E3").Select()
